$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Paragraph 4 currently exists empty (just the _GoBack bookmark).
# Fill it in with the first new bullet, then drop the _GoBack
# bookmark from here -- it will be recreated later, mid-word, in the
# final paragraph (that's where Word leaves it after the last edit).
# ------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertBefore("What are the assumptions of a Linear Regression?")

$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# Helper-free, explicit paragraph-by-paragraph construction so the
# list levels / text match the target exactly.
# ------------------------------------------------------------------

# 5: "No multi collinearity" (level 1) -- built so it ends up split
# into two runs, "No multi" + " collinearity", matching the source.
$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Item(5)
$p5.Range.ListFormat.ListLevelNumber = 2
$p5.Range.Text = "No multi collinearity"

$p5 = $d.Paragraphs.Item(5)
$firstWords = "No multi"
$splitPos = $p5.Range.Start + $firstWords.Length
$splitRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("TEMPSPLIT", $splitRange)
$d.Bookmarks.Item("TEMPSPLIT").Delete()

# 6: "The relationship between independent and dependent variable is linear" (level 1)
$p5 = $d.Paragraphs.Item(5)
$p5.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs.Item(6)
$p6.Range.ListFormat.ListLevelNumber = 2
$p6.Range.Text = "The relationship between independent and dependent variable is linear"

# 7: "The residuals of a linear regression follows normal distribution" (level 1)
$p6 = $d.Paragraphs.Item(6)
$p6.Range.InsertParagraphAfter()
$p7 = $d.Paragraphs.Item(7)
$p7.Range.ListFormat.ListLevelNumber = 2
$p7.Range.Text = "The residuals of a linear regression follows normal distribution"

# 8: "Gradient Descent is an optimization algorithm" (level 0)
$p7 = $d.Paragraphs.Item(7)
$p7.Range.InsertParagraphAfter()
$p8 = $d.Paragraphs.Item(8)
$p8.Range.ListFormat.ListLevelNumber = 1
$p8.Range.Text = "Gradient Descent is an optimization algorithm"

# 9: final paragraph (level 0), text split around the restored _GoBack bookmark.
$p8 = $d.Paragraphs.Item(8)
$p8.Range.InsertParagraphAfter()
$p9 = $d.Paragraphs.Item(9)
$p9.Range.ListFormat.ListLevelNumber = 1
$firstPart = "For linear regression, a mean squared error is a convex function and that is why Gradient Descent helps to optimize the cost function by subt"
$secondPart = "raction."
$p9.Range.Text = $firstPart + $secondPart

$p9 = $d.Paragraphs.Item(9)
$bmPos = $p9.Range.Start + $firstPart.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
